{"js": "// Replace the date line and each \"A\u00f7B=\" problem text with its updated value.\n// All old strings are unique within the document, so a simple search/replace\n// per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"2025-11-04 Tuesday\", \"2025-11-05 Wednesday\"],\n  [\"303\u00f73=\", \"528\u00f75=\"],\n  [\"946\u00f74=\", \"421\u00f77=\"],\n  [\"988\u00f73=\", \"840\u00f75=\"],\n  [\"773\u00f74=\", \"255\u00f73=\"],\n  [\"520\u00f73=\", \"541\u00f73=\"],\n  [\"981\u00f72=\", \"229\u00f79=\"],\n  [\"626\u00f75=\", \"688\u00f74=\"],\n  [\"800\u00f72=\", \"242\u00f73=\"],\n  [\"121\u00f77=\", \"700\u00f75=\"],\n  [\"467\u00f72=\", \"778\u00f78=\"],\n  [\"290\u00f77=\", \"702\u00f76=\"],\n  [\"103\u00f77=\", \"102\u00f78=\"],\n  [\"812\u00f73=\", \"289\u00f72=\"],\n  [\"812\u00f79=\", \"404\u00f75=\"],\n  [\"375\u00f75=\", \"947\u00f78=\"],\n  [\"178\u00f74=\", \"214\u00f78=\"],\n  [\"808\u00f74=\", \"491\u00f73=\"],\n  [\"444\u00f79=\", \"685\u00f75=\"],\n  [\"118\u00f75=\", \"960\u00f76=\"],\n  [\"575\u00f79=\", \"715\u00f76=\"],\n  [\"439\u00f75=\", \"868\u00f72=\"],\n  [\"401\u00f78=\", \"883\u00f74=\"],\n  [\"900\u00f76=\", \"438\u00f75=\"],\n  [\"997\u00f78=\", \"223\u00f76=\"],\n  [\"989\u00f75=\", \"713\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each \"A\u00f7B=\" division problem to its new value.\n# Every old string is unique in the document, so Find/Replace (wdReplaceAll)\n# per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2025-11-04 Tuesday\"; New = \"2025-11-05 Wednesday\"},\n    @{Old = \"303\u00f73=\"; New = \"528\u00f75=\"},\n    @{Old = \"946\u00f74=\"; New = \"421\u00f77=\"},\n    @{Old = \"988\u00f73=\"; New = \"840\u00f75=\"},\n    @{Old = \"773\u00f74=\"; New = \"255\u00f73=\"},\n    @{Old = \"520\u00f73=\"; New = \"541\u00f73=\"},\n    @{Old = \"981\u00f72=\"; New = \"229\u00f79=\"},\n    @{Old = \"626\u00f75=\"; New = \"688\u00f74=\"},\n    @{Old = \"800\u00f72=\"; New = \"242\u00f73=\"},\n    @{Old = \"121\u00f77=\"; New = \"700\u00f75=\"},\n    @{Old = \"467\u00f72=\"; New = \"778\u00f78=\"},\n    @{Old = \"290\u00f77=\"; New = \"702\u00f76=\"},\n    @{Old = \"103\u00f77=\"; New = \"102\u00f78=\"},\n    @{Old = \"812\u00f73=\"; New = \"289\u00f72=\"},\n    @{Old = \"812\u00f79=\"; New = \"404\u00f75=\"},\n    @{Old = \"375\u00f75=\"; New = \"947\u00f78=\"},\n    @{Old = \"178\u00f74=\"; New = \"214\u00f78=\"},\n    @{Old = \"808\u00f74=\"; New = \"491\u00f73=\"},\n    @{Old = \"444\u00f79=\"; New = \"685\u00f75=\"},\n    @{Old = \"118\u00f75=\"; New = \"960\u00f76=\"},\n    @{Old = \"575\u00f79=\"; New = \"715\u00f76=\"},\n    @{Old = \"439\u00f75=\"; New = \"868\u00f72=\"},\n    @{Old = \"401\u00f78=\"; New = \"883\u00f74=\"},\n    @{Old = \"900\u00f76=\"; New = \"438\u00f75=\"},\n    @{Old = \"997\u00f78=\"; New = \"223\u00f76=\"},\n    @{Old = \"989\u00f75=\"; New = \"713\u00f74=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$pair.New, [ref]2)\n}\n"}
